$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.638.84"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "3.447.42"
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.32"
$ws.Range("E5").Value = "  +1.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.49"
$ws.Range("E6").Value = "  +9.44%  "
$ws.Range("D7").Value = "3.448.61"
$ws.Range("E7").Value = "  +2.46%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +1.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.79"
$ws.Range("E10").Value = "  +2.00%  "
$ws.Range("E11").Value = "  +3.49%  "
$ws.Range("E12").Value = "  +1.58%  "
$ws.Range("D13").Value = "4.034.66"
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.92"
$ws.Range("E14").Value = "  +8.13%  "
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("E16").Value = "  +2.48%  "
$ws.Range("D17").Value = "3.447.00"
$ws.Range("E17").Value = "  +2.42%  "
$ws.Range("D18").Value = "61.726.44"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("E19").Value = "  +8.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.41"
$ws.Range("E20").Value = "  +3.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.50"
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.91"
$ws.Range("E22").Value = "  +4.05%  "
$ws.Range("E23").Value = "  +2.76%  "
$ws.Range("D24").Value = "3.588.64"
$ws.Range("E24").Value = "  +2.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.11"
$ws.Range("E25").Value = "  +3.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("E29").Value = "  +3.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.73"
$ws.Range("E30").Value = "  +3.84%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.54"
$ws.Range("E31").Value = "  -13.00%  "
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("E33").Value = "  +1.71%  "
$ws.Range("E34").Value = "  +1.50%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "24.09"
$ws.Range("E36").Value = "  +1.69%  "
$ws.Range("D37").Value = "3.475.12"
$ws.Range("E37").Value = "  +2.35%  "
$ws.Range("E38").Value = "  +2.76%  "
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("E42").Value = "  +3.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "27.20"
$ws.Range("E43").Value = "  +12.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.792"
$ws.Range("E44").Value = "  +2.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.50"
$ws.Range("E45").Value = "  +2.48%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "42.39"
$ws.Range("E47").Value = "  +2.04%  "
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("D49").Value = "2.602.41"
$ws.Range("E49").Value = "  +6.04%  "
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("E51").Value = "  +2.49%  "
